$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Header row (row 1), starting at column B like the "data" sheet
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"
$ws.Range("B1:G1").Font.Bold = $true
$ws.Range("B1:G1").HorizontalAlignment = -4108
$ws.Range("B1:G1").VerticalAlignment = -4160

# Data row (row 2)
$ws.Range("A2").Value = 0
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4160

$ws.Range("B2").Value = "Clefting"
$ws.Range("C2").Value = 81
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2.50"
$ws.Range("E2").Value = "2021-09-28T10:52:57.978494Z"
$ws.Range("F2").Value = "2021-10-05 14:19:36.573513"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/81/?format=json"

# Update the timestamps in the "data" sheet's time_taken column (F) to reflect
# the new run time recorded in the commit's diff.
$newTimes = @(
"2021-10-05 14:19:36.576986","2021-10-05 14:19:36.576995","2021-10-05 14:19:36.576999","2021-10-05 14:19:36.577001","2021-10-05 14:19:36.577004","2021-10-05 14:19:36.577007","2021-10-05 14:19:36.577010","2021-10-05 14:19:36.577012","2021-10-05 14:19:36.577015","2021-10-05 14:19:36.577018",
"2021-10-05 14:19:36.577021","2021-10-05 14:19:36.577023","2021-10-05 14:19:36.577026","2021-10-05 14:19:36.577029","2021-10-05 14:19:36.577032","2021-10-05 14:19:36.577034","2021-10-05 14:19:36.577037","2021-10-05 14:19:36.577040","2021-10-05 14:19:36.577043","2021-10-05 14:19:36.577046",
"2021-10-05 14:19:36.577049","2021-10-05 14:19:36.577051","2021-10-05 14:19:36.577054","2021-10-05 14:19:36.577057","2021-10-05 14:19:36.577059","2021-10-05 14:19:36.577062","2021-10-05 14:19:36.577065","2021-10-05 14:19:36.577067","2021-10-05 14:19:36.577070","2021-10-05 14:19:36.577073",
"2021-10-05 14:19:36.577075","2021-10-05 14:19:36.577078","2021-10-05 14:19:36.577081","2021-10-05 14:19:36.577083","2021-10-05 14:19:36.577086","2021-10-05 14:19:36.577089","2021-10-05 14:19:36.577091","2021-10-05 14:19:36.577094","2021-10-05 14:19:36.577096","2021-10-05 14:19:36.577099",
"2021-10-05 14:19:36.577102","2021-10-05 14:19:36.577105","2021-10-05 14:19:36.577108","2021-10-05 14:19:36.577110","2021-10-05 14:19:36.577112","2021-10-05 14:19:36.577115","2021-10-05 14:19:36.577118","2021-10-05 14:19:36.577120","2021-10-05 14:19:36.577123","2021-10-05 14:19:36.577126",
"2021-10-05 14:19:36.577129","2021-10-05 14:19:36.577131","2021-10-05 14:19:36.577134","2021-10-05 14:19:36.577137","2021-10-05 14:19:36.577140","2021-10-05 14:19:36.577143","2021-10-05 14:19:36.577145","2021-10-05 14:19:36.577148","2021-10-05 14:19:36.577151","2021-10-05 14:19:36.577153",
"2021-10-05 14:19:36.577156","2021-10-05 14:19:36.577158","2021-10-05 14:19:36.577161","2021-10-05 14:19:36.577163","2021-10-05 14:19:36.577167","2021-10-05 14:19:36.577170","2021-10-05 14:19:36.577173","2021-10-05 14:19:36.577175","2021-10-05 14:19:36.577178","2021-10-05 14:19:36.577181",
"2021-10-05 14:19:36.577184","2021-10-05 14:19:36.577187","2021-10-05 14:19:36.577189","2021-10-05 14:19:36.577192","2021-10-05 14:19:36.577195","2021-10-05 14:19:36.577198","2021-10-05 14:19:36.577202","2021-10-05 14:19:36.577206","2021-10-05 14:19:36.577208","2021-10-05 14:19:36.577211",
"2021-10-05 14:19:36.577214","2021-10-05 14:19:36.577216","2021-10-05 14:19:36.577219","2021-10-05 14:19:36.577222","2021-10-05 14:19:36.577224","2021-10-05 14:19:36.577227","2021-10-05 14:19:36.577230","2021-10-05 14:19:36.577233","2021-10-05 14:19:36.577236","2021-10-05 14:19:36.577239",
"2021-10-05 14:19:36.577241","2021-10-05 14:19:36.577244","2021-10-05 14:19:36.577248","2021-10-05 14:19:36.577251","2021-10-05 14:19:36.577254","2021-10-05 14:19:36.577256","2021-10-05 14:19:36.577259","2021-10-05 14:19:36.577262","2021-10-05 14:19:36.577264","2021-10-05 14:19:36.577267",
"2021-10-05 14:19:36.577270","2021-10-05 14:19:36.577272","2021-10-05 14:19:36.577275","2021-10-05 14:19:36.577277","2021-10-05 14:19:36.577280","2021-10-05 14:19:36.577283","2021-10-05 14:19:36.577285","2021-10-05 14:19:36.577288","2021-10-05 14:19:36.577293","2021-10-05 14:19:36.577296",
"2021-10-05 14:19:36.577299","2021-10-05 14:19:36.577302","2021-10-05 14:19:36.577305","2021-10-05 14:19:36.577307","2021-10-05 14:19:36.577310","2021-10-05 14:19:36.577312","2021-10-05 14:19:36.577315","2021-10-05 14:19:36.577318","2021-10-05 14:19:36.577320","2021-10-05 14:19:36.577323",
"2021-10-05 14:19:36.577325","2021-10-05 14:19:36.577328","2021-10-05 14:19:36.577331","2021-10-05 14:19:36.577333","2021-10-05 14:19:36.577336","2021-10-05 14:19:36.577338","2021-10-05 14:19:36.577341","2021-10-05 14:19:36.577343","2021-10-05 14:19:36.577348","2021-10-05 14:19:36.577351",
"2021-10-05 14:19:36.577354","2021-10-05 14:19:36.577356","2021-10-05 14:19:36.577363","2021-10-05 14:19:36.577366","2021-10-05 14:19:36.577369","2021-10-05 14:19:36.577372","2021-10-05 14:19:36.577374","2021-10-05 14:19:36.577377","2021-10-05 14:19:36.577380","2021-10-05 14:19:36.577382",
"2021-10-05 14:19:36.577385","2021-10-05 14:19:36.577388","2021-10-05 14:19:36.577390","2021-10-05 14:19:36.577393","2021-10-05 14:19:36.577396","2021-10-05 14:19:36.577398","2021-10-05 14:19:36.577401","2021-10-05 14:19:36.577403","2021-10-05 14:19:36.577406","2021-10-05 14:19:36.577409",
"2021-10-05 14:19:36.577412","2021-10-05 14:19:36.577414","2021-10-05 14:19:36.577417","2021-10-05 14:19:36.577419","2021-10-05 14:19:36.577422","2021-10-05 14:19:36.577425","2021-10-05 14:19:36.577427","2021-10-05 14:19:36.577430","2021-10-05 14:19:36.577433","2021-10-05 14:19:36.577435",
"2021-10-05 14:19:36.577438","2021-10-05 14:19:36.577441","2021-10-05 14:19:36.577443","2021-10-05 14:19:36.577446","2021-10-05 14:19:36.577448","2021-10-05 14:19:36.577451","2021-10-05 14:19:36.577454","2021-10-05 14:19:36.577456","2021-10-05 14:19:36.577459","2021-10-05 14:19:36.577461",
"2021-10-05 14:19:36.577464","2021-10-05 14:19:36.577466","2021-10-05 14:19:36.577470","2021-10-05 14:19:36.577474","2021-10-05 14:19:36.577476","2021-10-05 14:19:36.577479","2021-10-05 14:19:36.577481","2021-10-05 14:19:36.577484","2021-10-05 14:19:36.577486","2021-10-05 14:19:36.577489",
"2021-10-05 14:19:36.577492","2021-10-05 14:19:36.577494","2021-10-05 14:19:36.577497","2021-10-05 14:19:36.577499","2021-10-05 14:19:36.577502","2021-10-05 14:19:36.577505","2021-10-05 14:19:36.577507","2021-10-05 14:19:36.577510","2021-10-05 14:19:36.577513","2021-10-05 14:19:36.577515",
"2021-10-05 14:19:36.577518","2021-10-05 14:19:36.577520","2021-10-05 14:19:36.577523","2021-10-05 14:19:36.577526","2021-10-05 14:19:36.577529","2021-10-05 14:19:36.577531","2021-10-05 14:19:36.577534","2021-10-05 14:19:36.577536","2021-10-05 14:19:36.577539","2021-10-05 14:19:36.577542",
"2021-10-05 14:19:36.577545","2021-10-05 14:19:36.577547","2021-10-05 14:19:36.577550","2021-10-05 14:19:36.577552","2021-10-05 14:19:36.577555","2021-10-05 14:19:36.577558","2021-10-05 14:19:36.577560","2021-10-05 14:19:36.577562","2021-10-05 14:19:36.577565","2021-10-05 14:19:36.577568",
"2021-10-05 14:19:36.577570","2021-10-05 14:19:36.577573","2021-10-05 14:19:36.577575","2021-10-05 14:19:36.577578","2021-10-05 14:19:36.577581","2021-10-05 14:19:36.577583","2021-10-05 14:19:36.577586","2021-10-05 14:19:36.577588","2021-10-05 14:19:36.577591","2021-10-05 14:19:36.577594",
"2021-10-05 14:19:36.577596","2021-10-05 14:19:36.577599","2021-10-05 14:19:36.577601","2021-10-05 14:19:36.577604","2021-10-05 14:19:36.577607","2021-10-05 14:19:36.577609","2021-10-05 14:19:36.577612","2021-10-05 14:19:36.577615","2021-10-05 14:19:36.577617","2021-10-05 14:19:36.577620",
"2021-10-05 14:19:36.577623","2021-10-05 14:19:36.577625","2021-10-05 14:19:36.577629","2021-10-05 14:19:36.577632","2021-10-05 14:19:36.577635","2021-10-05 14:19:36.577638","2021-10-05 14:19:36.577641","2021-10-05 14:19:36.577643","2021-10-05 14:19:36.577646","2021-10-05 14:19:36.577648",
"2021-10-05 14:19:36.577651","2021-10-05 14:19:36.577653","2021-10-05 14:19:36.577656","2021-10-05 14:19:36.577659","2021-10-05 14:19:36.577661","2021-10-05 14:19:36.577664","2021-10-05 14:19:36.577667","2021-10-05 14:19:36.577669","2021-10-05 14:19:36.577672","2021-10-05 14:19:36.577674",
"2021-10-05 14:19:36.577677","2021-10-05 14:19:36.577679","2021-10-05 14:19:36.577682","2021-10-05 14:19:36.577685","2021-10-05 14:19:36.577687","2021-10-05 14:19:36.577690","2021-10-05 14:19:36.577692","2021-10-05 14:19:36.577695","2021-10-05 14:19:36.577697","2021-10-05 14:19:36.577700",
"2021-10-05 14:19:36.577703","2021-10-05 14:19:36.577705","2021-10-05 14:19:36.577708","2021-10-05 14:19:36.577710"
)

$data = $wb.Worksheets.Item("data")
for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 6).Value = $newTimes[$i]
}

Write-Output "done"
